$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "66.120.35"
$ws.Range("E2").Value2 = "  -0.12%  "
$ws.Range("D3").Value2 = "3.035.12"
$ws.Range("E3").Value2 = "  +0.85%  "
$ws.Range("E4").Value2 = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "575.61"
$ws.Range("E5").Value2 = "  -1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "168.43"
$ws.Range("E6").Value2 = "  +3.20%  "
$ws.Range("E7").Value2 = "  +0.07%  "
$ws.Range("D8").Value2 = "3.030.12"
$ws.Range("E8").Value2 = "  +0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "6.67"
$ws.Range("E10").Value2 = "  -0.23%  "
$ws.Range("E11").Value2 = "  -1.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.479"
$ws.Range("E12").Value2 = "  +5.00%  "
$ws.Range("E13").Value2 = "  -3.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "36.82"
$ws.Range("E14").Value2 = "  +6.30%  "
$ws.Range("E15").Value2 = "  -0.47%  "
$ws.Range("D16").Value2 = "66.131.22"
$ws.Range("E16").Value2 = "  -0.08%  "
$ws.Range("D17").Value2 = "3.547.68"
$ws.Range("E17").Value2 = "  +1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "7.21"
$ws.Range("E18").Value2 = "  +3.89%  "
$ws.Range("B19").Value2 = "Chainlink"
$ws.Range("C19").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "16.48"
$ws.Range("E19").Value2 = "  +18.26%  "
$ws.Range("B20").Value2 = "WrappedEther"
$ws.Range("C20").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value2 = "3.044.75"
$ws.Range("E20").Value2 = "  +1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "466.72"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.705"
$ws.Range("E22").Value2 = "  +2.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "7.46"
$ws.Range("E23").Value2 = "  +1.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "83.23"
$ws.Range("E24").Value2 = "  +1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "12.90"
$ws.Range("E25").Value2 = "  +4.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "2.27"
$ws.Range("E26").Value2 = "  -0.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "10.08"
$ws.Range("E27").Value2 = "  -4.44%  "
$ws.Range("E28").Value2 = "  +0.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "8.36"
$ws.Range("E29").Value2 = "  +2.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "2.45"
$ws.Range("E30").Value2 = "  +2.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "2.62"
$ws.Range("E31").Value2 = "  +0.46%  "
$ws.Range("E32").Value2 = "  +7.19%  "
$ws.Range("D33").Value2 = "0.0₃0992"
$ws.Range("E33").Value2 = "  -6.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "28.14"
$ws.Range("E34").Value2 = "  +3.28%  "
$ws.Range("E35").Value2 = "  +0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.993"
$ws.Range("E36").Value2 = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "5.82"
$ws.Range("E37").Value2 = "  +0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "48.64"
$ws.Range("E38").Value2 = "  +10.80%  "
$ws.Range("E39").Value2 = "  -6.93%  "
$ws.Range("E40").Value2 = "  -0.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.308"
$ws.Range("E41").Value2 = "  +0.15%  "
$ws.Range("E42").Value2 = "  -1.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "2.84"
$ws.Range("E43").Value2 = "  -6.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "8.60"
$ws.Range("E44").Value2 = "  +1.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.0358"
$ws.Range("E45").Value2 = "  -0.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "382.61"
$ws.Range("E46").Value2 = "  -3.59%  "
$ws.Range("D47").Value2 = "2.736.04"
$ws.Range("E47").Value2 = "  -2.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "134.30"
$ws.Range("E48").Value2 = "  +0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "24.81"
$ws.Range("E50").Value2 = "  +3.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "2.22"
$ws.Range("E51").Value2 = "  +3.33%  "
